# ---------------------------------------------------------------------------
# Edit: collapse the three "CORE COMPETENCIES" detail paragraphs into a
# single summary line, and append a new "TECHNICAL SKILLS" section (heading
# + three detail paragraphs) right before the closing "For a more
# detailed..." line.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument
$bullet = [char]0x2022

function Find-ParaIndex($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.StartsWith($needle)) {
            return $i
        }
    }
    throw "Paragraph starting with '$needle' not found"
}

# --- 1. Collapse the "CORE COMPETENCIES" detail paragraphs ----------------

$firstIdx = Find-ParaIndex $d "Data Visualization & Design: Interactive Dashboards"
$firstPara = $d.Paragraphs.Item($firstIdx)

# Replace the text of the first paragraph (leaving its own paragraph mark
# intact) with the new condensed summary line.
$bodyRange = $d.Range($firstPara.Range.Start, $firstPara.Range.End - 1)
$bodyRange.Text = "Data Visualization & Design $bullet Geospatial Analysis & Mapping $bullet Technical Visualization"

# Delete the next two paragraphs (and their marks) entirely, since their
# content has been folded into the line above.
$midIdx = Find-ParaIndex $d "Geospatial Analysis & Mapping: Spatial Analysis"
$midPara = $d.Paragraphs.Item($midIdx)
$endIdx = Find-ParaIndex $d "Technical Visualization: Programming: Python"
$endPara = $d.Paragraphs.Item($endIdx)

$deleteRange = $d.Range($midPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()

# --- 2. Insert the new "TECHNICAL SKILLS" section --------------------------
#
# All four new paragraphs are split off of the ("Normal"-styled) anchor
# paragraph one after another *before* any style is changed, so each new
# paragraph inherits the plain "Normal" style. Only afterwards do we promote
# the first of the four (the "TECHNICAL SKILLS" line) to Heading 2 -- this
# avoids the heading style leaking into the three detail paragraphs that
# follow it.

$anchorNeedle = "$bullet Led multi-million dollar research projects"
$anchorIdx = Find-ParaIndex $d $anchorNeedle
$anchorPara = $d.Paragraphs.Item($anchorIdx)

[void]$anchorPara.Range.InsertParagraphAfter()
$headingIdx = $anchorIdx + 1

[void]$d.Paragraphs.Item($headingIdx).Range.InsertParagraphAfter()
$idx1 = $headingIdx + 1

[void]$d.Paragraphs.Item($idx1).Range.InsertParagraphAfter()
$idx2 = $idx1 + 1

[void]$d.Paragraphs.Item($idx2).Range.InsertParagraphAfter()
$idx3 = $idx2 + 1

$d.Paragraphs.Item($headingIdx).Range.Text = "TECHNICAL SKILLS"
$d.Paragraphs.Item($idx1).Range.Text = "DATA VISUALIZATION & DESIGN Interactive Dashboards; Statistical Visualization; Geospatial Mapping; Choropleth Design"
$d.Paragraphs.Item($idx2).Range.Text = "GEOSPATIAL ANALYSIS & MAPPING Spatial Analysis; Mapping Technologies; Web Mapping; Spatial Data Processing"
$d.Paragraphs.Item($idx3).Range.Text = "TECHNICAL VISUALIZATION Programming; Database Integration; Web Technologies; Statistical Computing"

# Promote only the heading paragraph to Heading 2, last.
$d.Paragraphs.Item($headingIdx).Style = "Heading 2"

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
